# Commit message: "Changing the sign from + to - for the material recycled
# for each component"
#
# The workbook contains one worksheet per year (tabs "2000" .. "2100").
# Every worksheet shares the identical layout:
#   Row 1         : column headers (Generator Onshore / Offshore / Panel / Wires)
#   Rows 2-5, B:E : recycled-material quantities (rows = Nd, Dy, Cu, Si)
#
# The edit simply flips the sign of every non-zero quantity in the B2:E5
# block on every sheet (cells that are already 0 stay 0).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $rng = $ws.Range("B2:E5")
    foreach ($cell in $rng.Cells) {
        $v = $cell.Value2
        if ($v -ne 0) {
            $cell.Value2 = -$v
        }
    }
}
